# Insert a new "Fecha" observation row (row 110) for the weekly Mango series.
# This shifts the existing rows 110-192 down to 111-193 (Excel preserves all
# cell contents/styles automatically), and we then populate the new row 110
# with a duplicate of the entry that is now on row 111 (i.e. what used to be
# row 110), except for a new, more recent Fecha value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 110:192 down to 111:193, leaving row 110 blank.
$ws.Rows(110).Insert()

# Populate the new row 110 with the same record as the row below it (which
# now holds the data that used to live in row 110), but with an updated date.
$ws.Range("A110").Value = $ws.Range("A111").Value()
$ws.Range("B110").Value = $ws.Range("B111").Value()
$ws.Range("C110").Value = $ws.Range("C111").Value()
$ws.Range("D110").Value = 44673
$ws.Range("E110").Value = $ws.Range("E111").Value()
$ws.Range("F110").Value = $ws.Range("F111").Value()
$ws.Range("G110").Value = $ws.Range("G111").Value()
$ws.Range("H110").Value = $ws.Range("H111").Value()
$ws.Range("I110").Value = $ws.Range("I111").Value()
$ws.Range("J110").Value = $ws.Range("J111").Value()
$ws.Range("K110").Value = $ws.Range("K111").Value()
$ws.Range("L110").Value = $ws.Range("L111").Value()
$ws.Range("M110").Value = $ws.Range("M111").Value()
$ws.Range("N110").Value = $ws.Range("N111").Value()
$ws.Range("O110").Value = $ws.Range("O111").Value()
$ws.Range("P110").Value = $ws.Range("P111").Value()
$ws.Range("Q110").Value = $ws.Range("Q111").Value()
$ws.Range("R110").Value = $ws.Range("R111").Value()
$ws.Range("S110").Value = $ws.Range("S111").Value()
$ws.Range("T110").Value = $ws.Range("T111").Value()

# D110 carries the same date style (style index 2, yyyy-mm-dd-ish date number
# format) as the rest of the Fecha column.
$ws.Range("D110").Style = $ws.Range("D111").Style()
